$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - copy formatting from H1 (style index 1) and set labels
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J24
$data = @(
    @(6, 6),
    @(10, 10),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 9),
    @(9, 9),
    @(8, 9),
    @(6, 7),
    @(9, 9),
    @(7, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(7, 8),
    @(7, 8),
    @(9, 9),
    @(5, 6),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 9).Value = $data[$i][0]
    $ws.Cells.Item($r, 10).Value = $data[$i][1]
}

